$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Mónica Ávila" / "Padrino" (row 32) entirely, shifting
# everything below it up by one row.
$ws.Rows.Item(32).Delete()

# Fix the total formula to cover the now-shorter range (row shift already
# adjusts the SUM range automatically, but make sure it matches expectation).
$ws.Range("B35").Formula = "=SUM(B2:B34)"

# Restore the selection Excel leaves after this kind of edit.
$ws.Range("A8").Select()
